$d = $word.ActiveDocument

# Fix 1: "Load the file "CT.jpg" and display it" -> "...lungs.jpg..."
$d.Content.Find.Execute("CT.jpg", $false, $false, $false, $false, $false, $true, 1, $false, "lungs.jpg", 2)

# Fix 2: typo "lindow" -> "window"
$d.Content.Find.Execute("lindow", $false, $false, $false, $false, $false, $true, 1, $false, "window", 2)
